$p = $ppt.ActivePresentation
Write-Host "Slides count:" $p.Slides.Count
$s = $p.Slides.Item(3)
Write-Host "Comments count on slide 3:" $s.Comments.Count
